# ---------------------------------------------------------------------------
# "part 6 in prog" - add a new worksheet "part6" that redoes the GPU-CPU-SSD
# comparison but replacing the SSD leg with CXL, and tidy up number formats
# on part4 / part5 so everything uses plain integers instead of
# thousand-separated integers.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. part4: extend the "sanity check" column down through row 29 and switch
#    the thousands-separated number format to a plain integer format.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("part4")

$ws4.Range("B2:D6").NumberFormat   = "0"
$ws4.Range("B11:D15").NumberFormat = "0"
$ws4.Range("B32:D36").NumberFormat = "0"

$ws4.Range("E17").Formula = "=IF((B17+C17)=D17,TRUE,FALSE)"
$ws4.Range("E18").Formula = "=IF((B18+C18)=D18,TRUE,FALSE)"
$ws4.Range("E19").Formula = "=IF((B19+C19)=D19,TRUE,FALSE)"
$ws4.Range("E20").Formula = "=IF((B20+C20)=D20,TRUE,FALSE)"
$ws4.Range("E21").Formula = "=IF((B21+C21)=D21,TRUE,FALSE)"
$ws4.Range("E22").Formula = "=IF((B22+C22)=D22,TRUE,FALSE)"
$ws4.Range("E23").Formula = "=IF((B23+C23)=D23,TRUE,FALSE)"
$ws4.Range("E24").Formula = "=IF((B24+C24)=D24,TRUE,FALSE)"
$ws4.Range("E25").Formula = "=IF((B25+C25)=D25,TRUE,FALSE)"
$ws4.Range("E26").Formula = "=IF((B26+C26)=D26,TRUE,FALSE)"
$ws4.Range("E27").Formula = "=IF((B27+C27)=D27,TRUE,FALSE)"
$ws4.Range("E28").Formula = "=IF((B28+C28)=D28,TRUE,FALSE)"
$ws4.Range("E29").Formula = "=IF((B29+C29)=D29,TRUE,FALSE)"

$ws4.Range("A31").Select()

# ---------------------------------------------------------------------------
# 2. part5: switch number format to plain integer too (cosmetic match-up
#    with part4 / the new part6).
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("part5")
$ws5.Range("B1:D4").NumberFormat  = "0"
$ws5.Range("B9:D16").NumberFormat = "0"
$ws5.Range("B1").Select()

# ---------------------------------------------------------------------------
# 3. Add the new "part6" worksheet after "part5" and populate it.
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Add($null, $ws5)
$ws6.Name = "part6"

$ws6.Range("A1").Value = "BASED ON GPU-CPU-SSD, REPLACING SSD WITH CXL"

$ws6.Range("B2").Value = "Iter0"
$ws6.Range("C2").Value = "Iter1"
$ws6.Range("D2").Value = "Total"
$ws6.Range("E2").Value = "Sanity Check"
$ws6.Range("F2").Value = "Perf Ratio"

# -- Bert --------------------------------------------------------------
$ws6.Range("A3").Value = "Bert-base"
$ws6.Range("B3").Value = 131986612062
$ws6.Range("C3").Value = 112266689356
$ws6.Range("D3").Value = 244253301418

$ws6.Range("A4").Value = "Bert"
$ws6.Range("B4").Value = 126175516592
$ws6.Range("C4").Value = 112266689356
$ws6.Range("D4").Value = 238442205948
$ws6.Range("E4").Formula = "=IF(C4+B4=D4,TRUE,FALSE)"
$ws6.Range("F4").Formula = "=D3/D4"

# -- InceptionV3 ---------------------------------------------------------
$ws6.Range("A5").Value = "InceptionV3-base"
$ws6.Range("B5").Value = 142970763281
$ws6.Range("C5").Value = 139989684384
$ws6.Range("D5").Value = 282960447665

$ws6.Range("A6").Value = "InceptionV3"
$ws6.Range("B6").Value = 144700567589
$ws6.Range("C6").Value = 139989684384
$ws6.Range("D6").Value = 284690251973
$ws6.Range("E6").Formula = "=IF(C6+B6=D6,TRUE,FALSE)"
$ws6.Range("F6").Formula = "=D5/D6"

# -- ResNet152 -------------------------------------------------------------
$ws6.Range("A7").Value = "ResNet152-base"
$ws6.Range("B7").Value = 251655078132
$ws6.Range("C7").Value = 244319509089
$ws6.Range("D7").Value = 495974587221

$ws6.Range("A8").Value = "ResNet152"
$ws6.Range("B8").Value = 256372356123
$ws6.Range("C8").Value = 244319509089
$ws6.Range("D8").Value = 500691865212
$ws6.Range("E8").Formula = "=IF(C8+B8=D8,TRUE,FALSE)"
$ws6.Range("F8").Formula = "=D7/D8"

# -- SENet154 --------------------------------------------------------------
$ws6.Range("A9").Value = "SENet154-base"
$ws6.Range("B9").Value = 466228995826
$ws6.Range("C9").Value = 456431950066
$ws6.Range("D9").Value = 922660945892

$ws6.Range("A10").Value = "SENet154"
$ws6.Range("B10").Value = 471919199046
$ws6.Range("C10").Value = 456431950066
$ws6.Range("D10").Value = 928351149112
$ws6.Range("E10").Formula = "=IF(C10+B10=D10,TRUE,FALSE)"
$ws6.Range("F10").Formula = "=D9/D10"

# -- VIT ---------------------------------------------------------------
$ws6.Range("A11").Value = "VIT-base"
$ws6.Range("B11").Value = 25313940949
$ws6.Range("C11").Value = 21454177966
$ws6.Range("D11").Value = 46768118915

$ws6.Range("A12").Value = "VIT"
$ws6.Range("B12").Value = 23629607726
$ws6.Range("C12").Value = 21454177966
$ws6.Range("D12").Value = 45083785692
$ws6.Range("E12").Formula = "=IF(C12+B12=D12,TRUE,FALSE)"
$ws6.Range("F12").Formula = "=D11/D12"

# -- empty formatted rows (carried over layout below the table) -----------
$ws6.Range("B17:D17").Value = ""
$ws6.Range("B18:D22").Value = ""

# -- number formats / alignment -------------------------------------------
$ws6.Range("B2:D2").NumberFormat  = "0"
$ws6.Range("F2").NumberFormat     = "0"
$ws6.Range("B2:D2,F2").HorizontalAlignment = -4152  # xlRight

$ws6.Range("B3:D3,B5:D5,B7:D7,B9:D9,B11:D11").NumberFormat = "0"
$ws6.Range("B3:D3,B5:D5,B7:D7,B9:D9,B11:D11").HorizontalAlignment = -4152

$ws6.Range("B4:D4,B6:D6,B8:D8,B10:D10,B12:D12").NumberFormat = "0"

$ws6.Range("B17:D17").NumberFormat = "0"
$ws6.Range("B17:D17").HorizontalAlignment = -4152

$ws6.Range("B18:D22").NumberFormat = "0"
$ws6.Range("B18:D22").HorizontalAlignment = -4152

# -- column widths ----------------------------------------------------------
$ws6.Columns.Item("A").ColumnWidth = 22.61328125
$ws6.Columns.Item("B").ColumnWidth = 15.4609375
$ws6.Columns.Item("C").ColumnWidth = 13.53515625
$ws6.Columns.Item("D").ColumnWidth = 18.69140625
$ws6.Columns.Item("E").ColumnWidth = 18.53515625

$ws6.Range("A4").Select()
